# Adds row 13 to the "Artfynd" sheet, matching the appended observation
# record from the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 112205187
$ws.Range("B13").Value = 90662
$ws.Range("C13").Value = "Ovaliderad"
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 4363
$ws.Range("F13").Value = "Zontaggsvamp"
$ws.Range("G13").Value = "Hydnellum concrescens"
$ws.Range("H13").Value = "(Pers.) Banker"

# I13 holds the numeric-looking text "1" (not the number 1). Force text
# storage via NumberFormat, then drop back to the Normal style so no
# stray formatting is left behind.
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = "1"
$ws.Range("I13").Style = "Normal"

$ws.Range("J13").Value = "mycel"

# K13 / N13 / AF13 / AT13 / AY13 are present-but-empty text cells in the
# source row (<c t="inlineStr"/> with no content). A plain "" assignment
# clears/removes the cell instead of leaving an empty string behind, so
# use the quote-prefix trick (a lone "'" is Excel's marker for "empty
# text"), then strip the quote-prefix formatting it leaves behind.
$ws.Range("K13").Value = "'"
$ws.Range("K13").Style = "Normal"

$ws.Range("N13").Value = "'"
$ws.Range("N13").Style = "Normal"

$ws.Range("P13").Value = "Baggetorp, Nrk"
$ws.Range("Q13").Value = 478586.2573052422
$ws.Range("R13").Value = 6556136.916654737
$ws.Range("S13").Value = 10
$ws.Range("T13").Value = "Örebro"
$ws.Range("U13").Value = "Lekeberg"
$ws.Range("V13").Value = "Närke"
$ws.Range("W13").Value = "Kvistbro"

# Y13 / AA13 hold the literal date-like text "2023-09-19"; left alone,
# Excel's smart-typing recognises the string as a date and stores a
# serial number instead. Force text storage the same way as I13.
$ws.Range("Y13").NumberFormat = "@"
$ws.Range("Y13").Value = "2023-09-19"
$ws.Range("Y13").Style = "Normal"

$ws.Range("Z13").Value = "00:00"

$ws.Range("AA13").NumberFormat = "@"
$ws.Range("AA13").Value = "2023-09-19"
$ws.Range("AA13").Style = "Normal"

$ws.Range("AB13").Value = "00:00"

$ws.Range("AD13").Value = $false
$ws.Range("AE13").Value = $false

$ws.Range("AF13").Value = "'"
$ws.Range("AF13").Style = "Normal"

$ws.Range("AG13").Value = $false
$ws.Range("AH13").Value = "Blåbärsbarrskog"
$ws.Range("AI13").Value = "i yta bökad av vildsvin"

$ws.Range("AT13").Value = "'"
$ws.Range("AT13").Style = "Normal"

$ws.Range("AW13").Value = "Michael Andersson"
$ws.Range("AX13").Value = "Michael Andersson"

$ws.Range("AY13").Value = "'"
$ws.Range("AY13").Style = "Normal"
